$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.054.75'
$ws.Range('E2').Value = '  +1.71%  '

$ws.Range('D3').Value = '1.975.95'
$ws.Range('E3').Value = '  +4.69%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9950'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.52%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8104'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +71.82%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '252.98'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +3.68%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9931'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.69%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3429'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +18.47%  '

$ws.Range('E9').Value = '  +14.85%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06907'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +6.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8531'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +17.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08160'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +5.20%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '102.04'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +6.44%  '

$ws.Range('D14').Value = '1.968.45'
$ws.Range('E14').Value = '  +4.26%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.502'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +6.03%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '279.55'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -1.08%  '

$ws.Range('D17').Value = '31.008.16'
$ws.Range('E17').Value = '  +1.63%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.88'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +6.14%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007851'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +5.03%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.675'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +7.63%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.227.78'
$ws.Range('E21').Value = '  +4.27%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9952'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9966'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.30%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.754'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +6.70%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1656'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +70.97%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.620'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +5.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.86'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +0.50%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.54'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +3.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.189'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +15.57%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.560'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +6.11%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.358'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +1.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.537'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +6.04%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.330'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +4.30%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05116'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +5.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.219'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +8.22%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7389'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +6.44%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.753'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +1.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9936'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -0.56%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01983'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +5.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.893'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +2.38%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.575'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +6.00%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '78.42'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +4.50%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4650'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +8.87%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.068'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +4.32%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8471'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +2.32%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.60'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +3.16%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9940'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -0.60%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.975'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +3.33%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.457'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +7.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4260'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +8.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.23'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +2.93%  '
